$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the SVC row (was row 8); this shifts subsequent rows up
# and automatically removes "SVC" from the shared strings table.
$ws.Rows.Item(8).Delete()

# Rename NN_1 -> NN (now located at row 10, column A)
$ws.Cells.Item(10, 1).Value = "NN"

# Update numeric values for rows 4-12 (columns B..L) with the corrected statistics
# Row 4: LR
$ws.Cells.Item(4, 2).Value = 120.3474
$ws.Cells.Item(4, 3).Value = 19.0054
$ws.Cells.Item(4, 4).Value = 51.6475
$ws.Cells.Item(4, 5).Value = 5.3029
$ws.Cells.Item(4, 6).Value = 1.876
$ws.Cells.Item(4, 7).Value = 0.3848
$ws.Cells.Item(4, 8).Value = 0.819
$ws.Cells.Item(4, 9).Value = 0.1009
$ws.Cells.Item(4, 10).Value = 340.9474
$ws.Cells.Item(4, 11).Value = 74.7089
$ws.Cells.Item(4, 12).Value = 0.0148

# Row 5: Lasso
$ws.Cells.Item(5, 2).Value = 72.6348
$ws.Cells.Item(5, 3).Value = 7.7878
$ws.Cells.Item(5, 4).Value = 51.1205
$ws.Cells.Item(5, 5).Value = 4.7924
$ws.Cells.Item(5, 6).Value = 0.6314
$ws.Cells.Item(5, 7).Value = 0.005
$ws.Cells.Item(5, 8).Value = 0.584
$ws.Cells.Item(5, 9).Value = 0.005
$ws.Cells.Item(5, 10).Value = 80.94929999999999
$ws.Cells.Item(5, 11).Value = 2.3645
$ws.Cells.Item(5, 12).Value = 0.0155

# Row 6: Ridge
$ws.Cells.Item(6, 2).Value = 55.9408
$ws.Cells.Item(6, 3).Value = 5.7965
$ws.Cells.Item(6, 4).Value = 39.1098
$ws.Cells.Item(6, 5).Value = 3.5804
$ws.Cells.Item(6, 6).Value = 0.5062
$ws.Cells.Item(6, 7).Value = 0.0057
$ws.Cells.Item(6, 8).Value = 0.4628
$ws.Cells.Item(6, 9).Value = 0.0055
$ws.Cells.Item(6, 10).Value = 53.9989
$ws.Cells.Item(6, 11).Value = 1.8146
$ws.Cells.Item(6, 12).Value = 0.0152

# Row 7: ElasticNet
$ws.Cells.Item(7, 2).Value = 67.06
$ws.Cells.Item(7, 3).Value = 7.0605
$ws.Cells.Item(7, 4).Value = 46.519
$ws.Cells.Item(7, 5).Value = 4.4266
$ws.Cells.Item(7, 6).Value = 0.5815
$ws.Cells.Item(7, 7).Value = 0.0054
$ws.Cells.Item(7, 8).Value = 0.5251
$ws.Cells.Item(7, 9).Value = 0.0053
$ws.Cells.Item(7, 10).Value = 66.6511
$ws.Cells.Item(7, 11).Value = 2.3179
$ws.Cells.Item(7, 12).Value = 0.0152

# Row 8: SVR
$ws.Cells.Item(8, 2).Value = 34.4107
$ws.Cells.Item(8, 3).Value = 4.9783
$ws.Cells.Item(8, 4).Value = 20.0086
$ws.Cells.Item(8, 5).Value = 2.8058
$ws.Cells.Item(8, 6).Value = 0.3109
$ws.Cells.Item(8, 7).Value = 0.0052
$ws.Cells.Item(8, 8).Value = 0.2436
$ws.Cells.Item(8, 9).Value = 0.0037
$ws.Cells.Item(8, 10).Value = 26.9103
$ws.Cells.Item(8, 11).Value = 0.8249
$ws.Cells.Item(8, 12).Value = 0.0147

# Row 9: XGB
$ws.Cells.Item(9, 2).Value = 36.4319
$ws.Cells.Item(9, 3).Value = 4.0658
$ws.Cells.Item(9, 4).Value = 22.0222
$ws.Cells.Item(9, 5).Value = 1.9934
$ws.Cells.Item(9, 6).Value = 0.3877
$ws.Cells.Item(9, 7).Value = 0.0054
$ws.Cells.Item(9, 8).Value = 0.3173
$ws.Cells.Item(9, 9).Value = 0.004
$ws.Cells.Item(9, 10).Value = 29.6116
$ws.Cells.Item(9, 11).Value = 0.5973000000000001
$ws.Cells.Item(9, 12).Value = 0.076

# Row 10: NN
$ws.Cells.Item(10, 2).Value = 38.3938
$ws.Cells.Item(10, 3).Value = 4.594
$ws.Cells.Item(10, 4).Value = 25.3869
$ws.Cells.Item(10, 5).Value = 3.0841
$ws.Cells.Item(10, 6).Value = 0.3451
$ws.Cells.Item(10, 7).Value = 0.0052
$ws.Cells.Item(10, 8).Value = 0.294
$ws.Cells.Item(10, 9).Value = 0.0041
$ws.Cells.Item(10, 10).Value = 30.7188
$ws.Cells.Item(10, 11).Value = 1.0807
$ws.Cells.Item(10, 12).Value = 1.2017

# Row 11: GP
$ws.Cells.Item(11, 2).Value = 26.3299
$ws.Cells.Item(11, 3).Value = 5.2115
$ws.Cells.Item(11, 4).Value = 14.5611
$ws.Cells.Item(11, 5).Value = 2.0609
$ws.Cells.Item(11, 6).Value = 0.2442
$ws.Cells.Item(11, 7).Value = 0.0047
$ws.Cells.Item(11, 8).Value = 0.2004
$ws.Cells.Item(11, 9).Value = 0.0035
$ws.Cells.Item(11, 10).Value = 24.8684
$ws.Cells.Item(11, 11).Value = 0.7047
$ws.Cells.Item(11, 12).Value = 0.1974

# Row 12: Cat
$ws.Cells.Item(12, 2).Value = 29.738
$ws.Cells.Item(12, 3).Value = 3.4688
$ws.Cells.Item(12, 4).Value = 17.955
$ws.Cells.Item(12, 5).Value = 1.7338
$ws.Cells.Item(12, 6).Value = 0.3106
$ws.Cells.Item(12, 7).Value = 0.004
$ws.Cells.Item(12, 8).Value = 0.2556
$ws.Cells.Item(12, 9).Value = 0.0031
$ws.Cells.Item(12, 10).Value = 25.4918
$ws.Cells.Item(12, 11).Value = 0.5377999999999999
$ws.Cells.Item(12, 12).Value = 0.6296
